$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the new Date column cells as Text so date-like strings are not
# auto-converted into date serial numbers (matches existing inline-string cells).
$ws.Range("A634:A653").NumberFormat = "@"

$ws.Cells.Item(634, 1).Value = "2024-09-02"
$ws.Range("C634").Value = 650.9500122070312
$ws.Range("D634").Value = 1470.050048828125
$ws.Range("E634").Value = 608.5800170898438
$ws.Range("F634").Value = 1316.800048828125
$ws.Range("G634").Value = 811.2000122070312
$ws.Range("H634").Value = 17987.74047851562
$ws.Range("I634").Value = 0
$ws.Range("J634").Value = 193.5023999705204

$ws.Cells.Item(635, 1).Value = "2024-09-03"
$ws.Range("C635").Value = 640.0499877929688
$ws.Range("D635").Value = 1460.75
$ws.Range("E635").Value = 599.9400024414062
$ws.Range("F635").Value = 1341.949951171875
$ws.Range("G635").Value = 822.3499755859375
$ws.Range("H635").Value = 17977.66967773438
$ws.Range("I635").Value = -0.0005598702512568748
$ws.Range("J635").Value = 193.3940637332301

$ws.Cells.Item(636, 1).Value = "2024-09-04"
$ws.Range("C636").Value = 645.5999755859375
$ws.Range("D636").Value = 1475.300048828125
$ws.Range("E636").Value = 609
$ws.Range("F636").Value = 1327.75
$ws.Range("G636").Value = 824.2999877929688
$ws.Range("H636").Value = 18052.54992675781
$ws.Range("I636").Value = 0.004165181047695957
$ws.Range("J636").Value = 194.1995850222286

$ws.Cells.Item(637, 1).Value = "2024-09-05"
$ws.Range("C637").Value = 643.9000244140625
$ws.Range("D637").Value = 1457.699951171875
$ws.Range("E637").Value = 602.1799926757812
$ws.Range("F637").Value = 1254.800048828125
$ws.Range("G637").Value = 835.4000244140625
$ws.Range("H637").Value = 17792.94024658203
$ws.Range("I637").Value = -0.01438077619112318
$ws.Range("J637").Value = 191.406844253615

$ws.Cells.Item(638, 1).Value = "2024-09-06"
$ws.Range("C638").Value = 665.25
$ws.Range("D638").Value = 1443.449951171875
$ws.Range("E638").Value = 597.2999877929688
$ws.Range("F638").Value = 1256.849975585938
$ws.Range("G638").Value = 832.7000122070312
$ws.Range("H638").Value = 17880.34979248047
$ws.Range("I638").Value = 0.004912597057432855
$ws.Range("J638").Value = 192.3471489534678

$ws.Cells.Item(639, 1).Value = "2024-09-09"
$ws.Range("C639").Value = 676
$ws.Range("D639").Value = 1492.050048828125
$ws.Range("E639").Value = 610.3400268554688
$ws.Range("F639").Value = 1225.25
$ws.Range("G639").Value = 827.5999755859375
$ws.Range("H639").Value = 18025.32012939453
$ws.Range("I639").Value = 0.008107802061849448
$ws.Range("J639").Value = 193.9066615643436

$ws.Cells.Item(640, 1).Value = "2024-09-10"
$ws.Range("C640").Value = 680
$ws.Range("D640").Value = 1503.050048828125
$ws.Range("E640").Value = 608
$ws.Range("F640").Value = 1246
$ws.Range("G640").Value = 824.75
$ws.Range("H640").Value = 18130.15014648438
$ws.Range("I640").Value = 0.005815709032478913
$ws.Range("J640").Value = 195.0343662874612

$ws.Cells.Item(641, 1).Value = "2024-09-11"
$ws.Range("C641").Value = 680.4500122070312
$ws.Range("D641").Value = 1499.949951171875
$ws.Range("E641").Value = 627.6599731445312
$ws.Range("F641").Value = 1229
$ws.Range("G641").Value = 814
$ws.Range("H641").Value = 18088.97985839844
$ws.Range("I641").Value = -0.002270818926114677
$ws.Range("J641").Value = 194.5914785572528

$ws.Cells.Item(642, 1).Value = "2024-09-12"
$ws.Range("C642").Value = 686.0999755859375
$ws.Range("D642").Value = 1513.449951171875
$ws.Range("E642").Value = 645.5999755859375
$ws.Range("F642").Value = 1224.849975585938
$ws.Range("G642").Value = 809.7000122070312
$ws.Range("H642").Value = 18193.19958496094
$ws.Range("I642").Value = 0.00576150382046627
$ws.Range("J642").Value = 195.7126181043906

$ws.Cells.Item(643, 1).Value = "2024-09-13"
$ws.Range("C643").Value = 681.9500122070312
$ws.Range("D643").Value = 1491.300048828125
$ws.Range("E643").Value = 646.6500244140625
$ws.Range("F643").Value = 1229.300048828125
$ws.Range("G643").Value = 788.0499877929688
$ws.Range("H643").Value = 18027.60040283203
$ws.Range("I643").Value = -0.009102257211853799
$ws.Range("J643").Value = 193.9311915147991

$ws.Cells.Item(644, 1).Value = "2024-09-16"
$ws.Range("C644").Value = 695.2000122070312
$ws.Range("D644").Value = 1456.349975585938
$ws.Range("E644").Value = 621.0499877929688
$ws.Range("F644").Value = 1219.699951171875
$ws.Range("G644").Value = 751.9500122070312
$ws.Range("H644").Value = 17765.49987792969
$ws.Range("I644").Value = -0.01453884704817227
$ws.Range("J644").Value = 191.1116555834957

$ws.Cells.Item(645, 1).Value = "2024-09-17"
$ws.Range("C645").Value = 692
$ws.Range("D645").Value = 1459.400024414062
$ws.Range("E645").Value = 649.6500244140625
$ws.Range("F645").Value = 1222.949951171875
$ws.Range("G645").Value = 746.75
$ws.Range("H645").Value = 17827
$ws.Range("I645").Value = 0.003461772676980224
$ws.Range("J645").Value = 191.7732406910471

$ws.Cells.Item(646, 1).Value = "2024-09-18"
$ws.Range("C646").Value = 695.2999877929688
$ws.Range("D646").Value = 1432.150024414062
$ws.Range("E646").Value = 646.7000122070312
$ws.Range("F646").Value = 1224.550048828125
$ws.Range("G646").Value = 744.5999755859375
$ws.Range("H646").Value = 17755.70007324219
$ws.Range("I646").Value = -0.003999547133999692
$ws.Range("J646").Value = 191.0062345758634

$ws.Cells.Item(647, 1).Value = "2024-09-19"
$ws.Range("C647").Value = 697
$ws.Range("D647").Value = 1444.849975585938
$ws.Range("E647").Value = 652.1500244140625
$ws.Range("F647").Value = 1197.849975585938
$ws.Range("G647").Value = 747.2000122070312
$ws.Range("H647").Value = 17752.34997558594
$ws.Range("I647").Value = -0.0001886773060161447
$ws.Range("J647").Value = 190.9701960340913

$ws.Cells.Item(648, 1).Value = "2024-09-20"
$ws.Range("C648").Value = 709
$ws.Range("D648").Value = 1456.599975585938
$ws.Range("E648").Value = 654.4500122070312
$ws.Range("F648").Value = 1206.300048828125
$ws.Range("G648").Value = 747.5499877929688
$ws.Range("H648").Value = 17905.25006103516
$ws.Range("I648").Value = 0.008612949027001824
$ws.Range("J648").Value = 192.6150125982095

$ws.Cells.Item(649, 1).Value = "2024-09-23"
$ws.Range("C649").Value = 702.5
$ws.Range("D649").Value = 1449.300048828125
$ws.Range("E649").Value = 654.0999755859375
$ws.Range("F649").Value = 1190
$ws.Range("G649").Value = 763.75
$ws.Range("H649").Value = 17852.70007324219
$ws.Range("I649").Value = -0.00293489270542646
$ws.Range("J649").Value = 192.0497082027794

$ws.Cells.Item(650, 1).Value = "2024-09-24"
$ws.Range("C650").Value = 705.0999755859375
$ws.Range("D650").Value = 1446.349975585938
$ws.Range("E650").Value = 646.8499755859375
$ws.Range("F650").Value = 1194.699951171875
$ws.Range("G650").Value = 760.9500122070312
$ws.Range("H650").Value = 17843.19958496094
$ws.Range("I650").Value = -0.0005321597429113499
$ws.Range("J650").Value = 191.947507079436

$ws.Cells.Item(651, 1).Value = "2024-09-25"
$ws.Range("C651").Value = 689.2000122070312
$ws.Range("D651").Value = 1429.550048828125
$ws.Range("E651").Value = 633.2999877929688
$ws.Range("F651").Value = 1175.349975585938
$ws.Range("G651").Value = 742.5499877929688
$ws.Range("H651").Value = 17509.20007324219
$ws.Range("I651").Value = -0.0187185885652627
$ws.Range("J651").Value = 188.3545206682882

$ws.Cells.Item(652, 1).Value = "2024-09-26"
$ws.Range("C652").Value = 693.5999755859375
$ws.Range("D652").Value = 1422.300048828125
$ws.Range("E652").Value = 626.8499755859375
$ws.Range("F652").Value = 1165
$ws.Range("G652").Value = 742.25
$ws.Range("H652").Value = 17466.64990234375
$ws.Range("I652").Value = -0.002430160756656341
$ws.Range("J652").Value = 187.8967889038213

$ws.Cells.Item(653, 1).Value = "2024-09-27"
$ws.Range("C653").Value = 692.4500122070312
$ws.Range("D653").Value = 1388.650024414062
$ws.Range("E653").Value = 608.5499877929688
$ws.Range("F653").Value = 1165.550048828125
$ws.Range("G653").Value = 751.6500244140625
$ws.Range("H653").Value = 17342.00036621094
$ws.Range("I653").Value = -0.007136430673868747
$ws.Range("J653").Value = 186.5558764959666
